# Quarterly indexing esoteric bug-fix operation
#
# Column A holds the "quarter" date label for each GDP QoQ observation.
# Every date was stamped on the 1st of the quarter's first month (e.g.
# 1988-07-01 for Q3 1988). The fix re-stamps each label on the 15th of the
# month that FOLLOWS the quarter's first month instead (e.g. 1988-08-15).
# That is: new_date = DATE(YEAR(old_date), MONTH(old_date) + 1, 15)
#
# Column B (the QoQ growth values) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's day-0 epoch (serial date arithmetic uses the classic 1900 system,
# i.e. day 1 = 1900-01-01 with the well-known 1900-02-29 leap bug baked in,
# which is equivalent to counting whole days from 1899-12-30).
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1

# Row 1 is the header ("date" / "value"); data starts on row 2.
$dataStartRow = $firstRow + 1

for ($r = $dataStartRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2

    if ($null -eq $serial -or $serial -eq "") {
        continue
    }

    $oldDate   = $epoch.AddDays([double]$serial)
    $nextMonth = $oldDate.AddMonths(1)
    $newDate   = Get-Date -Year $nextMonth.Year -Month $nextMonth.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value2 = $newDate
}
